$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.568.32"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "2.882.11"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.87"
$ws.Range("E5").Value = "  -4.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.44"
$ws.Range("E6").Value = "  -3.37%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("D9").Value = "2.879.22"
$ws.Range("E9").Value = "  -2.34%  "

$ws.Range("E10").Value = "  -2.66%  "

$ws.Range("E11").Value = "  -2.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.428"
$ws.Range("E12").Value = "  -2.12%  "

$ws.Range("E13").Value = "  -1.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.55"
$ws.Range("E14").Value = "  -2.98%  "

$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "3.362.97"
$ws.Range("E16").Value = "  -2.17%  "

$ws.Range("D17").Value = "61.546.19"
$ws.Range("E17").Value = "  -2.06%  "

$ws.Range("D18").Value = "2.884.96"
$ws.Range("E18").Value = "  -2.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  -2.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.37"
$ws.Range("E20").Value = "  -2.21%  "

$ws.Range("E21").Value = "  -3.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.650"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.78"
$ws.Range("E23").Value = "  -3.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.77"
$ws.Range("E24").Value = "  -2.34%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.83"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.89"
$ws.Range("E27").Value = "  -12.29%  "

$ws.Range("E28").Value = "  -6.03%  "

$ws.Range("E29").Value = "  +5.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.02"
$ws.Range("E30").Value = "  -4.03%  "

$ws.Range("E31").Value = "  -4.77%  "

$ws.Range("E32").Value = "  -9.62%  "

$ws.Range("E33").Value = "  +0.05%  "

$ws.Range("E34").Value = "  -2.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.32"
$ws.Range("E35").Value = "  -3.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.953"
$ws.Range("E36").Value = "  -3.81%  "

$ws.Range("E37").Value = "  -4.69%  "

$ws.Range("E39").Value = "  -7.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.89"
$ws.Range("E40").Value = "  -6.32%  "

$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.14"
$ws.Range("E41").Value = "  -3.46%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.113"
$ws.Range("E42").Value = "  -3.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.33"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("E44").Value = "  -5.12%  "

$ws.Range("D45").Value = "2.675.58"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.49"
$ws.Range("E46").Value = "  -1.83%  "

$ws.Range("E47").Value = "  -1.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "341.95"
$ws.Range("E48").Value = "  -4.41%  "

$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.38"
$ws.Range("E51").Value = "  -5.65%  "
